# Update rows 37-40 on the active sheet to reflect the new product ordering
# (a row for "Intel i7-11700K Rocket Lake 3.60GHz" was moved up, shifting the
# other processor rows down/around it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Processador Intel i7-11700K Rocket Lake 3.60GHz", "R$ 2.795,31", "https://www.kabum.com.br/produto/497469/processador-intel-i7-11700k-rocket-lake-3-60ghz"),
    @("Processador Intel Core i7 11700 2.50GHz",          "R$ 2.414,61", "https://www.kabum.com.br/produto/496250/processador-intel-core-i7-11700-2-50ghz"),
    @("Processador Intel Core i7 10700KF 3.80GHz",        "R$ 2.213,91", "https://www.kabum.com.br/produto/497468/processador-intel-core-i7-10700kf-3-80ghz"),
    @("Processador Intel Core i7 10700K 3.80GHz",         "R$ 2.335,41", "https://www.kabum.com.br/produto/496626/processador-intel-core-i7-10700k-3-80ghz")
)

$startRow = 37
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

$wb.Save()
